$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(6, 20, 20, 20, 1, 2.5, 1, 0.5, 0, 15, 10, 0, 0, 0, 0.01, 0.01, -0.03940771786447463, 15.25981612461969, 7.256195618949778, 0, "eyescrunching"),
    @(6, 20, 20, 20, 1, 2.5, 1, 0.5, 0, 15, 10, 0, 0, 0, 0.01, 0.1,  -0.01612645253041128, 25.87237940441076, -6.713321467312299, 1, "jaw"),
    @(6, 20, 20, 20, 1, 2.5, 1, 0.5, 0, 15, 10, 0, 0, 0, 0.01, 0.01, -0.3533387901154966,  13.81645549551831, -0.3712374617917025, 2, "raisingeyebrows"),
    @(6, 20, 20, 20, 1, 2.5, 1, 0.5, 0, 15, 10, 0, 0, 0, 0.01, 0.01, -0.2536986705699699,  10.2209899937896,  7.658682778488473,  3, "movehat"),
    @(6, 20, 20, 20, 1, 2.5, 1, 0.5, 0, 15, 10, 0, 0, 0, 0.01, 0.01, -0.1598917497698968,  11.94148954035882, 11.62803806050134,  4, "movehead")
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
